$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.836.72'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '2.364.29'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '''0.672'
$ws.Range('E5').Value = '  -1.29%  '
$ws.Range('D6').Value = '''239.88'
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('D7').Value = '''74.23'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '''0.602'
$ws.Range('E9').Value = '  +1.12%  '
$ws.Range('E10').Value = '  +1.46%  '
$ws.Range('D11').Value = '''60.07'
$ws.Range('D12').Value = '''37.13'
$ws.Range('E12').Value = '  +14.57%  '
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').Value = '''0.928'
$ws.Range('E16').Value = '  +2.96%  '
$ws.Range('D17').Value = '2.416.24'
$ws.Range('E17').Value = '  +2.65%  '
$ws.Range('D18').Value = '43.825.39'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('E19').Value = '  +1.49%  '
$ws.Range('D20').Value = '''6.62'
$ws.Range('E20').Value = '  -1.96%  '
$ws.Range('D21').Value = '''77.34'
$ws.Range('E21').Value = '  +0.41%  '
$ws.Range('D22').Value = '''251.45'
$ws.Range('E22').Value = '  -2.93%  '
$ws.Range('E23').Value = '  +3.64%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').Value = '''1.90'
$ws.Range('E25').Value = '  -4.05%  '
$ws.Range('D26').Value = '''2.50'
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').Value = '''10.57'
$ws.Range('E27').Value = '  -2.01%  '
$ws.Range('D28').Value = '''2.30'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('D29').Value = '''22.37'
$ws.Range('E29').Value = '  -1.53%  '
$ws.Range('D30').Value = '''175.32'
$ws.Range('E30').Value = '  -0.47%  '
$ws.Range('D31').Value = '''0.129'
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').Value = '''5.44'
$ws.Range('E34').Value = '  -1.06%  '
$ws.Range('D35').Value = '''5.10'
$ws.Range('E35').Value = '  -2.69%  '
$ws.Range('D36').Value = '''3.80'
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('D37').Value = '''6.61'
$ws.Range('E37').Value = '  +4.74%  '
$ws.Range('E38').Value = '  +1.99%  '
$ws.Range('E39').Value = '  +0.89%  '
$ws.Range('D40').Value = '''5.67'
$ws.Range('E40').Value = '  +20.73%  '
$ws.Range('D41').Value = '''20.66'
$ws.Range('E41').Value = '  +8.71%  '
$ws.Range('D42').Value = '''65.26'
$ws.Range('E42').Value = '  +11.92%  '
$ws.Range('E43').Value = '  -2.73%  '
$ws.Range('D44').Value = '''9.08'
$ws.Range('E44').Value = '  +1.08%  '
$ws.Range('D45').Value = '''0.202'
$ws.Range('E45').Value = '  -1.37%  '
$ws.Range('D46').Value = '''2.54'
$ws.Range('E46').Value = '  +1.52%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('E48').Value = '  +0.22%  '
$ws.Range('E49').Value = '  -1.14%  '
$ws.Range('D50').Value = '''98.46'
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('D51').Value = '''4.41'
$ws.Range('E51').Value = '  +15.28%  '
